# MOD: - Stand ohne Ultraschall / Servo aus ZBW fixiert
$wb = $excel.ActiveWorkbook

$wsTabelle1    = $wb.Worksheets.Item("Tabelle1")
$wsFahrversuch = $wb.Worksheets.Item("Fahrversuche")

# --- Tabelle1: correct the fixed (no-ultrasonic) delay value used for the
#     "Geradeausfahren" calibration row (chart1 / trendline reference this
#     cell) ---
$wsTabelle1.Range("C12").Value = 6000

# --- Fahrversuche: add the new "mit Zielsystem" deviation / factor block ---
# Cells are written in the same order the values were authored so that any
# newly introduced shared strings land in the same table order.
$wsFahrversuch.Range("A12").Value = "Mit Zielsystem"
$wsFahrversuch.Range("C12").Value = "(delay Faktor 210"
$wsFahrversuch.Range("C13").Value = "aufwärts"
$wsFahrversuch.Range("D12").Value = "Angabe in cm"
$wsFahrversuch.Range("E12").Value = "min"
$wsFahrversuch.Range("F12").Value = "max"
$wsFahrversuch.Range("E11").Value = "tatsächlich gefahren:"
$wsFahrversuch.Range("G12").Value = "Abweichung max"
$wsFahrversuch.Range("H12").Value = "Abweichung min"
$wsFahrversuch.Range("I12").Value = "Abweichung mittel"
$wsFahrversuch.Range("J12").Value = "Faktor"

$wsFahrversuch.Range("D13").Value = 40
$wsFahrversuch.Range("E13").Value = 26
$wsFahrversuch.Range("F13").Value = 30
$wsFahrversuch.Range("G13").Formula = "=D13-E13"
$wsFahrversuch.Range("H13").Formula = "=D13-F13"
$wsFahrversuch.Range("I13").Formula = "=(G13+H13)/2"
$wsFahrversuch.Range("J13").Formula = "=I13/D13"

$wsFahrversuch.Range("I15").Formula = "=E13+E13*J13"

# Column widths for the newly filled-in columns (auto-fit to content;
# values chosen so the stored sheet width lands on/near the real Excel
# auto-fit widths: A=14, D=20.43, E=19.57, F=4.71, G=16.29, H=15.86, I=18,
# J=6.57 characters).
$wsFahrversuch.Columns.Item(1).ColumnWidth = 13.166666666666666
$wsFahrversuch.Columns.Item(4).ColumnWidth = 19.666666666666668
$wsFahrversuch.Columns.Item(5).ColumnWidth = 18.666666666666668
$wsFahrversuch.Columns.Item(6).ColumnWidth = 3.8333333333333335
$wsFahrversuch.Columns.Item(7).ColumnWidth = 15.499999999999998
$wsFahrversuch.Columns.Item(8).ColumnWidth = 15.0
$wsFahrversuch.Columns.Item(9).ColumnWidth = 17.166666666666668
$wsFahrversuch.Columns.Item(10).ColumnWidth = 5.666666666666667

# --- View state: the workbook was saved with "Fahrversuche" as the active
#     sheet/tab, Tabelle1's selection parked on C12 and Fahrversuche's
#     selection on the newly added summary formula I15 ---
$wsTabelle1.Range("C12").Select() | Out-Null
$wsFahrversuch.Activate() | Out-Null
$wsFahrversuch.Range("I15").Select() | Out-Null
